# "Generate Report for Archive"
# Refresh the localization-status report:
#   - Status moves from "Ready for handoff" to "In Translation" (Overview + per-locale sheets)
#   - A new handoff-name timestamp token is recorded for the latest handoff
#   - The now-stale 2017-11-04 handoff-name token is gone, replaced by the 2017-11-17 one
#   - The "Status"/locale-status columns got a bit narrower after the refresh

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newHandoffName = "LocaleLowerCaseTest_HT_OL#Test1#20171117T080113"

# --- zh-cn and de-de detail sheets -----------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column C = "Status" for rows 2-5
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C4").Value = $newStatus
    $ws.Range("C5").Value = $newStatus

    # Column I = "Lastest Handoff Name" for rows 2-5
    $ws.Range("I2").Value = $newHandoffName
    $ws.Range("I3").Value = $newHandoffName
    $ws.Range("I4").Value = $newHandoffName
    $ws.Range("I5").Value = $newHandoffName

    # Column C width narrowed as part of the refreshed report layout
    $ws.Columns("C").ColumnWidth = 12.5
}

# --- Overview sheet ----------------------------------------------------------
# Columns E (zh-cn) and F (de-de) mirror the same "Status" value, and the
# matching width change.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus
$overview.Range("E5").Value = $newStatus
$overview.Range("F5").Value = $newStatus

$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
